# Add 'Stored in desiccator' as a new storage method (Closes #45)
# Also brings this file's shared enumerations (preparation_medium,
# storage_medium) up to date with 'HPMC-PVP' (already present in the
# sample-section workbook per a previous commit), and refreshes the
# template's pav:createdOn timestamp.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. storage_method: insert "Stored in desiccator" after "Unknown"
#    (row 6), before "Incubated at 37 degrees celsius".
# ------------------------------------------------------------------
$wsStorageMethod = $wb.Worksheets.Item("storage_method")
$wsStorageMethod.Rows.Item(6).Insert()
$wsStorageMethod.Cells.Item(6, 1).Value = "Stored in desiccator"
$wsStorageMethod.Cells.Item(6, 2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000387"

# ------------------------------------------------------------------
# 2. preparation_medium: insert "HPMC-PVP" after "Modified Davidson's
#    Fixative" (row 18), before "Inflated (Agarose)".
# ------------------------------------------------------------------
$wsPrepMedium = $wb.Worksheets.Item("preparation_medium")
$wsPrepMedium.Rows.Item(19).Insert()
$wsPrepMedium.Cells.Item(19, 1).Value = "HPMC-PVP"
$wsPrepMedium.Cells.Item(19, 2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000386"

# ------------------------------------------------------------------
# 3. storage_medium: insert "HPMC-PVP" after "Formic acid in water"
#    (row 6), before "DMSO (no serum)".
# ------------------------------------------------------------------
$wsStorageMedium = $wb.Worksheets.Item("storage_medium")
$wsStorageMedium.Rows.Item(7).Insert()
$wsStorageMedium.Cells.Item(7, 1).Value = "HPMC-PVP"
$wsStorageMedium.Cells.Item(7, 2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000386"

# ------------------------------------------------------------------
# 4. tissue_weight_unit: reorder to ug, g, mg, kg (was ug, mg, kg, g).
# ------------------------------------------------------------------
$wsTissueWeightUnit = $wb.Worksheets.Item("tissue_weight_unit")
$wsTissueWeightUnit.Cells.Item(2, 1).Value = "g"
$wsTissueWeightUnit.Cells.Item(2, 2).Value = "http://purl.obolibrary.org/obo/UO_0000021"
$wsTissueWeightUnit.Cells.Item(3, 1).Value = "mg"
$wsTissueWeightUnit.Cells.Item(3, 2).Value = "http://purl.obolibrary.org/obo/UO_0000022"
$wsTissueWeightUnit.Cells.Item(4, 1).Value = "kg"
$wsTissueWeightUnit.Cells.Item(4, 2).Value = "http://purl.obolibrary.org/obo/UO_0000009"

# ------------------------------------------------------------------
# 5. is_suspension_enriched: reorder to Yes, No (was No, Yes).
# ------------------------------------------------------------------
$wsIsEnriched = $wb.Worksheets.Item("is_suspension_enriched")
$wsIsEnriched.Cells.Item(1, 1).Value = "Yes"
$wsIsEnriched.Cells.Item(2, 1).Value = "No"

# ------------------------------------------------------------------
# 6. .metadata: bump pav:createdOn.
# ------------------------------------------------------------------
$wsMetadata = $wb.Worksheets.Item(".metadata")
$wsMetadata.Cells.Item(2, 3).Value = "2024-06-28T15:48:18-07:00"

# ------------------------------------------------------------------
# 7. Update the dropdown (data validation) ranges on the main sheet so
#    they cover the newly inserted rows.
# ------------------------------------------------------------------
$wsMain = $wb.Worksheets.Item("Sample Suspension")
$wsMain.Range("I2:I1001").Validation.Formula1 = "'preparation_medium'!`$A`$1:`$A`$31"
$wsMain.Range("M2:M1001").Validation.Formula1 = "'storage_medium'!`$A`$1:`$A`$22"
$wsMain.Range("N2:N1001").Validation.Formula1 = "'storage_method'!`$A`$1:`$A`$12"
